$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B17 no longer carries the extra "applyFill" style (revert to default style / no fill)
$ws.Range("B17").Interior.Pattern = -4142

# C17 value changes from 123456 to 12345
$ws.Range("C17").Value = 12345

# Update the active selection to G13
$ws.Range("G13").Select()
